$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tracker rows appended 2025-09-20 (Excel serial date 45920), continuing
# the existing pattern of one row per goal.
$goals = @(
    @{ Id = "G2"; Name = "Workout" },
    @{ Id = "G3"; Name = "Eat Healthy" },
    @{ Id = "G4"; Name = "Read Book" },
    @{ Id = "G5"; Name = "Investment Plan" },
    @{ Id = "G6"; Name = "Spend 10 Hours without phone" }
)

$startRow = 67
$date = 45920
$progress = 0.8786625992724292
$percentage = 0
$change = -0.01

for ($i = 0; $i -lt $goals.Count; $i++) {
    $row = $startRow + $i
    $goal = $goals[$i]

    $ws.Cells.Item($row, 1).Value = $goal.Id
    $ws.Cells.Item($row, 2).Value = $goal.Name
    $ws.Cells.Item($row, 3).Value = $date
    $ws.Cells.Item($row, 3).NumberFormat = $ws.Cells.Item($row - 1, 3).NumberFormat
    $ws.Cells.Item($row, 4).Value = $progress
    $ws.Cells.Item($row, 5).Value = $percentage
    $ws.Cells.Item($row, 6).Value = $change
}
